$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 17567.584
$ws.Range("J32").Value = 17701.334
$ws.Range("L32").Value = 17701.334
$ws.Range("N32").Value = -18353.334
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H69").Value = 9319.388999999999
$ws.Range("I69").Value = 11750
$ws.Range("J69").Value = 9015.5625
$ws.Range("K69").Value = 35250
$ws.Range("L69").Value = 27046.6875
$ws.Range("M69").Value = -34376
$ws.Range("N69").Value = -28794.6875
$ws.Range("H72").Value = 9319.388999999999
$ws.Range("I72").Value = 11750
$ws.Range("J72").Value = 9015.5625
$ws.Range("K72").Value = 105750
$ws.Range("L72").Value = 81140.0625
$ws.Range("M72").Value = -101382
$ws.Range("N72").Value = -89876.0625
$ws.Range("H106").Value = 9431.556
$ws.Range("I106").Value = 3012.8
$ws.Range("K106").Value = 3012.8
$ws.Range("M106").Value = -2381.8
$ws.Range("H132").Value = 2415.2144
$ws.Range("I132").Value = 2317.75
$ws.Range("K132").Value = 6953.25
$ws.Range("M132").Value = -4423.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3999.5
$ws.Range("I63").Value = 3999.5
$ws.Range("K63").Value = 3999.5
$ws.Range("M63").Value = -3313.5
$ws.Range("H66").Value = 3999.5
$ws.Range("I66").Value = 3999.5
$ws.Range("K66").Value = 19997.5
$ws.Range("M66").Value = -16565.5
$ws.Range("H124").Value = 34147
$ws.Range("J124").Value = 34147
$ws.Range("L124").Value = 34147
$ws.Range("N124").Value = -43967
$ws.Range("H132").Value = 7206.1924
$ws.Range("I132").Value = 5572.2
$ws.Range("K132").Value = 16716.6
$ws.Range("M132").Value = -14186.6
$ws.Range("H135").Value = 56345
$ws.Range("J135").Value = 56345
$ws.Range("L135").Value = 56345
$ws.Range("N135").Value = -66485

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2240.6667
$ws.Range("I99").Value = 2361
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 2361
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -863
$ws.Range("N99").Value = -4996
$ws.Range("H103").Value = 46900.6
$ws.Range("J103").Value = 46900.6
$ws.Range("L103").Value = 46900.6
$ws.Range("N103").Value = -49244.6
$ws.Range("H134").Value = 4132
$ws.Range("I134").Value = 1970.8462
$ws.Range("J134").Value = 11155.75
$ws.Range("K134").Value = 5912.5386
$ws.Range("L134").Value = 33467.25
$ws.Range("M134").Value = -3377.5386
$ws.Range("N134").Value = -38537.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 18602
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 18602
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 18602
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -19850
$ws.Range("H65").Value = 18602
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 18602
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 93010
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -99250
$ws.Range("H99").Value = 3289
$ws.Range("I99").Value = 2170.6667
$ws.Range("J99").Value = 3960
$ws.Range("K99").Value = 2170.6667
$ws.Range("L99").Value = 3960
$ws.Range("M99").Value = -672.6667000000002
$ws.Range("N99").Value = -6956
$ws.Range("H108").Value = 64999.5
$ws.Range("J108").Value = 64999.5
$ws.Range("L108").Value = 64999.5
$ws.Range("N108").Value = -72679.5
$ws.Range("H126").Value = 3289
$ws.Range("I126").Value = 2170.6667
$ws.Range("J126").Value = 3960
$ws.Range("K126").Value = 6512.000100000001
$ws.Range("L126").Value = 11880
$ws.Range("M126").Value = -4042.000100000001
$ws.Range("N126").Value = -16820
$ws.Range("H132").Value = 4294.5557
$ws.Range("I132").Value = 4061.25
$ws.Range("J132").Value = 4868.846
$ws.Range("K132").Value = 12183.75
$ws.Range("L132").Value = 14606.538
$ws.Range("M132").Value = -9653.75
$ws.Range("N132").Value = -19666.538
$ws.Range("H141").Value = 204327.5
$ws.Range("J141").Value = 237193
$ws.Range("L141").Value = 237193
$ws.Range("N141").Value = -247553

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1676.7693
$ws.Range("I34").Value = 1512.375
$ws.Range("J34").Value = 1939.8
$ws.Range("K34").Value = 4537.125
$ws.Range("L34").Value = 5819.4
$ws.Range("M34").Value = -4453.125
$ws.Range("N34").Value = -5987.4
$ws.Range("H98").Value = 276.125
$ws.Range("J98").Value = 243.33333
$ws.Range("L98").Value = 729.99999
$ws.Range("N98").Value = -3725.99999
$ws.Range("H122").Value = 1684.5555
$ws.Range("I122").Value = 217.5
$ws.Range("J122").Value = 2103.7144
$ws.Range("K122").Value = 1957.5
$ws.Range("L122").Value = 18933.4296
$ws.Range("M122").Value = 492.5
$ws.Range("N122").Value = -23833.4296

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 488.90475
$ws.Range("I97").Value = 398.58823
$ws.Range("K97").Value = 398.58823
$ws.Range("M97").Value = 97.41176999999999
$ws.Range("H132").Value = 53880
$ws.Range("I132").Value = 64221.707
$ws.Range("J132").Value = 9927.75
$ws.Range("K132").Value = 192665.121
$ws.Range("L132").Value = 29783.25
$ws.Range("M132").Value = -190135.121
$ws.Range("N132").Value = -34843.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 17290.25
$ws.Range("J43").Value = 15621
$ws.Range("L43").Value = 15621
$ws.Range("N43").Value = -16007
$ws.Range("H46").Value = 4583.5
$ws.Range("I46").Value = 4101
$ws.Range("J46").Value = 4680
$ws.Range("K46").Value = 4101
$ws.Range("L46").Value = 4680
$ws.Range("M46").Value = -3913
$ws.Range("N46").Value = -5056
$ws.Range("H93").Value = 1390.6
$ws.Range("I93").Value = 1412.2916
$ws.Range("K93").Value = 1412.2916
$ws.Range("M93").Value = -164.2916
$ws.Range("H100").Value = 12958.857
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2040.8
$ws.Range("I100").Value = 1539.8
$ws.Range("J100").Value = 2541.8
$ws.Range("K100").Value = 3079.6
$ws.Range("L100").Value = 5083.6
$ws.Range("M100").Value = -2538.6
$ws.Range("N100").Value = -6165.6
